$d = $word.ActiveDocument
$d.Content.Find.Execute("as one of the 20 patches around patch A, and patch B.", $true, $false, $false, $false, $false,
                         $true, 1, $false, "as one of the 20 patches around patch A, and as one of the 20 patches around patch B.", 2)
